# Fruta / hortaliza, semanal
# Re-shuffle the weekly Fecha / Volumen / Precio mínimo / Precio máximo /
# Precio promedio ponderado / Precio $/Kg values among the existing data
# rows (2,3,6,7,8,9,10). Rows 4, 5 and 11 are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2  (was row 9)
$ws.Range("D2").Value = 44449
$ws.Range("J2").Value = 1300
$ws.Range("K2").Value = 900
$ws.Range("L2").Value = 950
$ws.Range("M2").Value = 925
$ws.Range("P2").Value = 925

# Row 3  (was row 6)
$ws.Range("D3").Value = 44284
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 800
$ws.Range("L3").Value = 850
$ws.Range("M3").Value = 825
$ws.Range("P3").Value = 825

# Row 6  (was row 8)
$ws.Range("D6").Value = 44341
$ws.Range("J6").Value = 1300
$ws.Range("K6").Value = 900
$ws.Range("L6").Value = 1000
$ws.Range("M6").Value = 950
$ws.Range("P6").Value = 950

# Row 7  (was row 10)
$ws.Range("D7").Value = 44442
$ws.Range("J7").Value = 1250
$ws.Range("K7").Value = 850
$ws.Range("L7").Value = 900
$ws.Range("M7").Value = 875
$ws.Range("P7").Value = 875

# Row 8  (was row 3)
$ws.Range("D8").Value = 44453
$ws.Range("J8").Value = 1000
$ws.Range("K8").Value = 800
$ws.Range("L8").Value = 900
$ws.Range("M8").Value = 850
$ws.Range("P8").Value = 850

# Row 9  (was row 2)
$ws.Range("D9").Value = 44243
$ws.Range("J9").Value = 1200
$ws.Range("K9").Value = 1200
$ws.Range("L9").Value = 1300
$ws.Range("M9").Value = 1250
$ws.Range("P9").Value = 1250

# Row 10 (was row 7)
$ws.Range("D10").Value = 44291
$ws.Range("J10").Value = 1000
$ws.Range("K10").Value = 1000
$ws.Range("L10").Value = 1200
$ws.Range("M10").Value = 1100
$ws.Range("P10").Value = 1100
